# Fixed material number handling
#
# Input!B3  (Pro number)  20583536      -> 20595735
# Input!B4  (Instrument SN) "A01612"    -> "A01672"
# Input!B5  (ICA SN)        "APXCAS2134009" -> 12341 (now a plain number)
#
# Template_printout pulls these values via formulas, so they recompute
# automatically once the Input sheet is updated.

$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("Input")

$wsInput.Range("B3").Value = 20595735
$wsInput.Range("B4").Value = "A01672"
$wsInput.Range("B5").Value = 12341

# Input was the active/selected sheet in the saved workbook.
$wsInput.Activate() | Out-Null
$wsInput.Range("B5").Select() | Out-Null
